$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2129.3809
$ws.Range("I33").Value = 1885
$ws.Range("K33").Value = 1885
$ws.Range("M33").Value = -1656

$ws.Range("H40").Value = 4094.4285
$ws.Range("I40").Value = 3883.3333
$ws.Range("K40").Value = 3883.3333
$ws.Range("M40").Value = -3708.3333

$ws.Range("H74").Value = 3195.5625
$ws.Range("I74").Value = 3086.5386
$ws.Range("K74").Value = 3086.5386
$ws.Range("M74").Value = -2150.5386

$ws.Range("H76").Value = 7673.8823
$ws.Range("I76").Value = 9036.375
$ws.Range("J76").Value = 6462.778
$ws.Range("K76").Value = 9036.375
$ws.Range("L76").Value = 6462.778
$ws.Range("M76").Value = -8721.375
$ws.Range("N76").Value = -7092.778

$ws.Range("H77").Value = 3195.5625
$ws.Range("I77").Value = 3086.5386
$ws.Range("K77").Value = 15432.693
$ws.Range("M77").Value = -10752.693

$ws.Range("H79").Value = 7673.8823
$ws.Range("I79").Value = 9036.375
$ws.Range("J79").Value = 6462.778
$ws.Range("K79").Value = 9036.375
$ws.Range("L79").Value = 6462.778
$ws.Range("M79").Value = -7944.375
$ws.Range("N79").Value = -8646.778

$ws.Range("H138").Value = 2310.886
$ws.Range("J138").Value = 2491.4614
$ws.Range("L138").Value = 7474.3842
$ws.Range("N138").Value = -17754.3842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 533
$ws.Range("I5").Value = 180.5
$ws.Range("J5").Value = 650.5
$ws.Range("K5").Value = 180.5
$ws.Range("L5").Value = 650.5
$ws.Range("M5").Value = -68.5
$ws.Range("N5").Value = -874.5

$ws.Range("H63").Value = 4949.7896
$ws.Range("I63").Value = 2366.818
$ws.Range("K63").Value = 2366.818
$ws.Range("M63").Value = -1680.818

$ws.Range("H66").Value = 4949.7896
$ws.Range("I66").Value = 2366.818
$ws.Range("K66").Value = 11834.09
$ws.Range("M66").Value = -8402.09

$ws.Range("H92").Value = 56472.25
$ws.Range("J92").Value = 56472.25
$ws.Range("L92").Value = 56472.25
$ws.Range("N92").Value = -61464.25

$ws.Range("H97").Value = 911.8570999999999
$ws.Range("I97").Value = 834.1111
$ws.Range("K97").Value = 834.1111
$ws.Range("M97").Value = -338.1111

$ws.Range("H110").Value = 2650.75
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 2650.75
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 2650.75
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = -6740.75

$ws.Range("H134").Value = 100000
$ws.Range("J134").Value = 100000
$ws.Range("L134").Value = 100000
$ws.Range("N134").Value = -110140

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 533
$ws.Range("I4").Value = 180.5
$ws.Range("J4").Value = 650.5
$ws.Range("K4").Value = 180.5
$ws.Range("L4").Value = 650.5
$ws.Range("M4").Value = -65.5
$ws.Range("N4").Value = -880.5

$ws.Range("H22").Value = 1250.5
$ws.Range("I22").Value = 1250.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1250.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1077.5
$ws.Range("N22").ClearContents()

$ws.Range("H62").Value = 122000
$ws.Range("J62").Value = 122000
$ws.Range("L62").Value = 122000
$ws.Range("N62").Value = -123372

$ws.Range("H65").Value = 122000
$ws.Range("J65").Value = 122000
$ws.Range("L65").Value = 366000
$ws.Range("N65").Value = -372864

$ws.Range("H94").Value = 582.86664
$ws.Range("I94").Value = 393.58334
$ws.Range("J94").Value = 1340
$ws.Range("K94").Value = 393.58334
$ws.Range("L94").Value = 1340
$ws.Range("M94").Value = 57.41665999999998
$ws.Range("N94").Value = -2242

$ws.Range("H105").Value = 2612.652
$ws.Range("I105").Value = 1465
$ws.Range("J105").Value = 2854.2632
$ws.Range("K105").Value = 1465
$ws.Range("L105").Value = 2854.2632
$ws.Range("M105").Value = 282
$ws.Range("N105").Value = -6348.263199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 98990
$ws.Range("J92").Value = 98990
$ws.Range("L92").Value = 98990
$ws.Range("N92").Value = -103982

$ws.Range("H107").Value = 9868.799999999999
$ws.Range("I107").Value = 12448
$ws.Range("J107").Value = 6000
$ws.Range("K107").Value = 12448
$ws.Range("L107").Value = 6000
$ws.Range("M107").Value = -10528
$ws.Range("N107").Value = -9840

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 132.97437
$ws.Range("I2").Value = 76.5
$ws.Range("J2").Value = 152.44827
$ws.Range("K2").Value = 459
$ws.Range("L2").Value = 914.6896200000001
$ws.Range("M2").Value = -346
$ws.Range("N2").Value = -1140.68962

$ws.Range("H38").Value = 104.8421
$ws.Range("I38").Value = 111.666664
$ws.Range("J38").Value = 93.14286
$ws.Range("K38").Value = 334.999992
$ws.Range("L38").Value = 279.42858
$ws.Range("M38").Value = 12.00000799999998
$ws.Range("N38").Value = -973.42858

$ws.Range("H97").Value = 8929293
$ws.Range("J97").Value = 2000
$ws.Range("L97").Value = 6000
$ws.Range("N97").Value = -6992

$ws.Range("H102").Value = 5000
$ws.Range("I102").Value = 2000
$ws.Range("K102").Value = 6000
$ws.Range("M102").Value = -3566

$ws.Range("H131").Value = 209757.81
$ws.Range("I131").Value = 1263075
$ws.Range("J131").Value = 18245.592
$ws.Range("K131").Value = 3789225
$ws.Range("L131").Value = 54736.776
$ws.Range("M131").Value = -3784185
$ws.Range("N131").Value = -64816.776

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 29990
$ws.Range("I34").Value = 29990
$ws.Range("K34").Value = 29990
$ws.Range("M34").Value = -29722

$ws.Range("H62").Value = 115000
$ws.Range("J62").Value = 115000
$ws.Range("L62").Value = 115000
$ws.Range("N62").Value = -116372

$ws.Range("H65").Value = 115000
$ws.Range("J65").Value = 115000
$ws.Range("L65").Value = 345000
$ws.Range("N65").Value = -351864

$ws.Range("H70").Value = 25165.666
$ws.Range("J70").Value = 8665
$ws.Range("L70").Value = 8665
$ws.Range("N70").Value = -9205

$ws.Range("H73").Value = 25165.666
$ws.Range("J73").Value = 8665
$ws.Range("L73").Value = 8665
$ws.Range("N73").Value = -10537

$ws.Range("H76").Value = 29990
$ws.Range("I76").Value = 29990
$ws.Range("K76").Value = 29990
$ws.Range("M76").Value = -29675

$ws.Range("H79").Value = 29990
$ws.Range("I79").Value = 29990
$ws.Range("K79").Value = 29990
$ws.Range("M79").Value = -28898

$ws.Range("H80").Value = 3613.5881
$ws.Range("I80").Value = 3347.3845
$ws.Range("J80").Value = 4478.75
$ws.Range("K80").Value = 3347.3845
$ws.Range("L80").Value = 4478.75
$ws.Range("M80").Value = -2349.3845
$ws.Range("N80").Value = -6474.75

$ws.Range("H83").Value = 3613.5881
$ws.Range("I83").Value = 3347.3845
$ws.Range("J83").Value = 4478.75
$ws.Range("K83").Value = 16736.9225
$ws.Range("L83").Value = 22393.75
$ws.Range("M83").Value = -11744.9225
$ws.Range("N83").Value = -32377.75

$ws.Range("H97").Value = 1345
$ws.Range("I97").Value = 1195.5385
$ws.Range("J97").Value = 1992.6666
$ws.Range("K97").Value = 1195.5385
$ws.Range("L97").Value = 1992.6666
$ws.Range("M97").Value = -699.5385000000001
$ws.Range("N97").Value = -2984.6666

$ws.Range("H105").Value = 39333
$ws.Range("J105").Value = 39333
$ws.Range("L105").Value = 39333
$ws.Range("N105").Value = -46321

$ws.Range("H134").Value = 150000
$ws.Range("J134").Value = 150000
$ws.Range("L134").Value = 450000
$ws.Range("N134").Value = -455070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8632.678
$ws.Range("I7").Value = 4115.2
$ws.Range("K7").Value = 4115.2
$ws.Range("M7").Value = -4003.2

$ws.Range("H46").Value = 3164.6155
$ws.Range("J46").Value = 4041.4285
$ws.Range("L46").Value = 4041.4285
$ws.Range("N46").Value = -4417.4285

$ws.Range("H126").Value = 8632.678
$ws.Range("I126").Value = 4115.2
$ws.Range("K126").Value = 12345.6
$ws.Range("M126").Value = -9875.599999999999

$ws.Range("H132").Value = 919358.4399999999
$ws.Range("I132").Value = 112921.78
$ws.Range("K132").Value = 338765.34
$ws.Range("M132").Value = -336235.34

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 100001050
$ws.Range("I107").Value = 166667550
$ws.Range("J107").Value = 1296
$ws.Range("K107").Value = 500002650
$ws.Range("L107").Value = 3888
$ws.Range("M107").Value = -500000730
$ws.Range("N107").Value = -7728

$ws.Range("H136").Value = 7309.385
$ws.Range("J136").Value = 6561.75
$ws.Range("L136").Value = 19685.25
$ws.Range("N136").Value = -24785.25
